$d = $word.ActiveDocument

# Locate the run that ends with " experiment" (end of the last paragraph) so the
# two new runs (a line break, then a line break + "one more time") get inserted
# right after it, inside the same paragraph.
$find = $d.Content
$find.Find.Execute("experiment", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$find.Collapse(0)  # wdCollapseEnd

$br = [string][char]11

# --- First new run: a single manual line break ---
$insertPos1 = $find.End
$r1 = $d.Range($insertPos1, $insertPos1)
$r1.Select()
$word.Selection.TypeText($br)
$word.Selection.LanguageID = "en-US"

# --- Second new run: another manual line break followed by the new text ---
$insertPos2 = $d.Content.End - 1
$r2 = $d.Range($insertPos2, $insertPos2)
$r2.Select()
$word.Selection.TypeText($br + "one more time")
$word.Selection.LanguageID = "en-US"
